$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1242.8823
$ws.Range("I100").Value = 1075
$ws.Range("J100").Value = 1550.6666
$ws.Range("K100").Value = 1075
$ws.Range("L100").Value = 1550.6666
$ws.Range("M100").Value = -534
$ws.Range("N100").Value = -2632.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("N25").ClearContents()

$ws.Range("H26").Value = 15351.75
$ws.Range("J26").Value = 19866.666
$ws.Range("L26").Value = 19866.666
$ws.Range("N26").Value = -20526.666

$ws.Range("H64").Value = 30000
$ws.Range("J64").Value = 30000
$ws.Range("L64").Value = 30000
$ws.Range("N64").Value = -30496

$ws.Range("H67").Value = 30000
$ws.Range("J67").Value = 30000
$ws.Range("L67").Value = 30000
$ws.Range("N67").Value = -31716

$ws.Range("H104").Value = 37707
$ws.Range("J104").Value = 37707
$ws.Range("L104").Value = 37707
$ws.Range("N104").Value = -44695

$ws.Range("H109").Value = 44988
$ws.Range("J109").Value = 44988
$ws.Range("L109").Value = 44988
$ws.Range("N109").Value = -47762

$ws.Range("H117").Value = 45643.285
$ws.Range("J117").Value = 45643.285
$ws.Range("L117").Value = 45643.285
$ws.Range("N117").Value = -54821.285

$ws.Range("H119").Value = 54807
$ws.Range("J119").Value = 54807
$ws.Range("L119").Value = 54807
$ws.Range("N119").Value = -64483

$ws.Range("H132").Value = 17858862
$ws.Range("I132").Value = 33334494
$ws.Range("K132").Value = 100003482
$ws.Range("M132").Value = -100000952

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 26620
$ws.Range("J19").Value = 32900
$ws.Range("L19").Value = 32900
$ws.Range("N19").Value = -33246

$ws.Range("H62").Value = 30000
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 30000
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

$ws.Range("H117").Value = 43497.332
$ws.Range("J117").Value = 43497.332
$ws.Range("L117").Value = 43497.332
$ws.Range("N117").Value = -52675.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 4592.5
$ws.Range("I36").Value = 4956.6665
$ws.Range("K36").Value = 4956.6665
$ws.Range("M36").Value = -4568.6665

$ws.Range("H38").Value = 25000
$ws.Range("I38").Value = 26000
$ws.Range("K38").Value = 26000
$ws.Range("M38").Value = -25623

$ws.Range("H40").Value = 4592.5
$ws.Range("I40").Value = 4956.6665
$ws.Range("K40").Value = 4956.6665
$ws.Range("M40").Value = -4796.6665

$ws.Range("H46").Value = 25000
$ws.Range("I46").Value = 26000
$ws.Range("K46").Value = 26000
$ws.Range("M46").Value = -25789

$ws.Range("H106").Value = 37146.25
$ws.Range("J106").Value = 37146.25
$ws.Range("L106").Value = 37146.25
$ws.Range("N106").Value = -39670.25

$ws.Range("H116").Value = 47668.5
$ws.Range("J116").Value = 47668.5
$ws.Range("L116").Value = 47668.5
$ws.Range("N116").Value = -56846.5

$ws.Range("H117").Value = 20000
$ws.Range("J117").Value = 20000
$ws.Range("L117").Value = 20000
$ws.Range("N117").Value = -29178

$ws.Range("H118").Value = 44710
$ws.Range("J118").Value = 44710
$ws.Range("L118").Value = 44710
$ws.Range("N118").Value = -48024

$ws.Range("H119").Value = 48504.332
$ws.Range("J119").Value = 48504.332
$ws.Range("L119").Value = 48504.332
$ws.Range("N119").Value = -58180.332

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 4043.0571
$ws.Range("I5").Value = 20459.4
$ws.Range("K5").Value = 61378.2
$ws.Range("M5").Value = -61266.2

$ws.Range("H9").Value = 166667680
$ws.Range("I9").Value = 200000700
$ws.Range("J9").Value = 2600
$ws.Range("K9").Value = 600002100
$ws.Range("L9").Value = 7800
$ws.Range("M9").Value = -600001876
$ws.Range("N9").Value = -8248

$ws.Range("H82").Value = 3671

$ws.Range("H85").Value = 3671

$ws.Range("H121").Value = 282215.2
$ws.Range("I121").Value = 282.5
$ws.Range("J121").Value = 443319.56
$ws.Range("K121").Value = 847.5
$ws.Range("L121").Value = 1329958.68
$ws.Range("M121").Value = 462.5
$ws.Range("N121").Value = -1332578.68

$ws.Range("H135").Value = 4043.0571
$ws.Range("I135").Value = 20459.4
$ws.Range("K135").Value = 184134.6
$ws.Range("M135").Value = -181599.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 2857.1428
$ws.Range("J12").Value = 5000
$ws.Range("L12").Value = 5000
$ws.Range("N12").Value = -5280

$ws.Range("H46").Value = 24952
$ws.Range("J46").Value = 24952
$ws.Range("L46").Value = 24952
$ws.Range("N46").Value = -25264

$ws.Range("H55").Value = 3950
$ws.Range("I55").Value = 3950
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 3950
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -3623
$ws.Range("N55").ClearContents()

$ws.Range("H105").Value = 39301
$ws.Range("J105").Value = 39301
$ws.Range("L105").Value = 39301
$ws.Range("N105").Value = -46289

$ws.Range("H107").Value = 589.6786
$ws.Range("I107").Value = 366.5
$ws.Range("J107").Value = 991.4
$ws.Range("K107").Value = 366.5
$ws.Range("L107").Value = 991.4
$ws.Range("M107").Value = 1553.5
$ws.Range("N107").Value = -4831.4

$ws.Range("H118").Value = 38086.668
$ws.Range("J118").Value = 38086.668
$ws.Range("L118").Value = 38086.668
$ws.Range("N118").Value = -41400.668

$ws.Range("H132").Value = 5324.36
$ws.Range("I132").Value = 2200.9167
$ws.Range("J132").Value = 8207.538
$ws.Range("K132").Value = 6602.750100000001
$ws.Range("L132").Value = 24622.614
$ws.Range("M132").Value = -4072.750100000001
$ws.Range("N132").Value = -29682.614

$ws.Range("H134").Value = 30448.334
$ws.Range("J134").Value = 30448.334
$ws.Range("L134").Value = 91345.00199999999
$ws.Range("N134").Value = -96415.00199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2891.1924
$ws.Range("I7").Value = 2479.125
$ws.Range("J7").Value = 3550.5
$ws.Range("K7").Value = 2479.125
$ws.Range("L7").Value = 3550.5
$ws.Range("M7").Value = -2367.125
$ws.Range("N7").Value = -3774.5

$ws.Range("H11").Value = 1000
$ws.Range("J11").Value = 1000
$ws.Range("L11").Value = 1000
$ws.Range("N11").Value = -1280

$ws.Range("H24").Value = 3603
$ws.Range("I24").Value = 4206
$ws.Range("J24").Value = 3000
$ws.Range("K24").Value = 4206
$ws.Range("L24").Value = 3000
$ws.Range("M24").Value = -3863
$ws.Range("N24").Value = -3686

$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()

$ws.Range("H108").Value = 48622
$ws.Range("J108").Value = 48622
$ws.Range("L108").Value = 48622
$ws.Range("N108").Value = -56302

$ws.Range("H119").Value = 47408
$ws.Range("J119").Value = 47408
$ws.Range("L119").Value = 47408
$ws.Range("N119").Value = -57084

$ws.Range("H120").Value = 46773.332
$ws.Range("J120").Value = 46773.332
$ws.Range("L120").Value = 46773.332
$ws.Range("N120").Value = -56449.332

$ws.Range("H126").Value = 2891.1924
$ws.Range("I126").Value = 2479.125
$ws.Range("J126").Value = 3550.5
$ws.Range("K126").Value = 7437.375
$ws.Range("L126").Value = 10651.5
$ws.Range("M126").Value = -4967.375
$ws.Range("N126").Value = -15591.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 45935.25
$ws.Range("J16").Value = 45935.25
$ws.Range("L16").Value = 45935.25
$ws.Range("N16").Value = -46519.25

$ws.Range("H104").Value = 37521
$ws.Range("J104").Value = 37521
$ws.Range("L104").Value = 37521
$ws.Range("N104").Value = -44509

$ws.Range("H105").Value = 43455
$ws.Range("J105").Value = 43455
$ws.Range("L105").Value = 43455
$ws.Range("N105").Value = -50443

$ws.Range("H114").Value = 34798.5
$ws.Range("I114").Value = 30000
$ws.Range("J114").Value = 36398
$ws.Range("K114").Value = 30000
$ws.Range("L114").Value = 36398
$ws.Range("M114").Value = -25661
$ws.Range("N114").Value = -45076

$ws.Range("H119").Value = 48694
$ws.Range("J119").Value = 48694
$ws.Range("L119").Value = 48694
$ws.Range("N119").Value = -58370

$ws.Range("H126").Value = 1133045.9
$ws.Range("I126").Value = 1338417.8
$ws.Range("J126").Value = 3501
$ws.Range("K126").Value = 4015253.4
$ws.Range("L126").Value = 10503
$ws.Range("M126").Value = -4012783.4
$ws.Range("N126").Value = -15443
